$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before D, shifting existing D:K data right to F:M.
# (New quarters Dec-2018 and Sep-2018 are being added to the front of the report,
#  pushing the historical quarters over by two columns.)
$ws.Range("D:E").Insert()

# Copy number/date formatting from column F (the shifted former column D)
# into the two new D:E columns so they inherit the correct per-row style.
# (Done per contiguous data block so the blank section-header rows 5, 6, 37 and 79
#  are not given stray formatted cells in D:E.)
$ws.Range("F7:F36").Copy()
$ws.Range("D7:E36").PasteSpecial(-4122)
$ws.Range("F38:F78").Copy()
$ws.Range("D38:E78").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Populate new column D (period ending 2018-12-31) and E (2018-09-30) values.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 272200
$ws.Range("E8").Value = 295100
$ws.Range("D9").Value = 115100
$ws.Range("E9").Value = 120400
$ws.Range("D10").Value = 157100
$ws.Range("E10").Value = 174700
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 6400
$ws.Range("E15").Value = 8100
$ws.Range("D17").Value = 227800
$ws.Range("E17").Value = 235600
$ws.Range("D18").Value = 44400
$ws.Range("E18").Value = 59500
$ws.Range("D20").Value = 17400
$ws.Range("E20").Value = 1700
$ws.Range("D21").Value = 70800
$ws.Range("E21").Value = 69300
$ws.Range("D22").Value = 1600
$ws.Range("E22").Value = 1600
$ws.Range("D23").Value = 60200
$ws.Range("E23").Value = 59600
$ws.Range("D24").Value = 19500
$ws.Range("E24").Value = 14100
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 40700
$ws.Range("E26").Value = 45500
$ws.Range("D27").Value = 41100
$ws.Range("E27").Value = 45300
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 5400
$ws.Range("E29").Value = 1000
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -17400
$ws.Range("E32").Value = -1700
$ws.Range("D33").Value = 46500
$ws.Range("E33").Value = 46300
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 46500
$ws.Range("E35").Value = 46300
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 232000
$ws.Range("E41").Value = 270500
$ws.Range("D42").Value = 617100
$ws.Range("E42").Value = 588400
$ws.Range("D43").Value = 169600
$ws.Range("E43").Value = 144100
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 86700
$ws.Range("E45").Value = 56900
$ws.Range("D46").Value = 1105500
$ws.Range("E46").Value = 1060000
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 63400
$ws.Range("E48").Value = 69300
$ws.Range("D49").Value = 145900
$ws.Range("E49").Value = 145900
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 29300
$ws.Range("E52").Value = 17600
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1344100
$ws.Range("E54").Value = 1292800
$ws.Range("D57").Value = 26300
$ws.Range("E57").Value = 40500
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = "NA"
$ws.Range("D59").Value = 311900
$ws.Range("E59").Value = 222200
$ws.Range("D60").Value = 338100
$ws.Range("E60").Value = 262600
$ws.Range("D61").Value = 94900
$ws.Range("E61").Value = 94800
$ws.Range("D62").Value = 16200
$ws.Range("E62").Value = 25700
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 460600
$ws.Range("E66").Value = 399300
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 1198400
$ws.Range("E72").Value = 1170800
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 883500
$ws.Range("E76").Value = 893500
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 46500
$ws.Range("E81").Value = 46300
$ws.Range("D83").Value = 9000
$ws.Range("E83").Value = 8100
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 93300
$ws.Range("E89").Value = 90700
$ws.Range("D91").Value = -700
$ws.Range("E91").Value = -700
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -32200
$ws.Range("E94").Value = -7200
$ws.Range("D96").Value = -19700
$ws.Range("E96").Value = -20100
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -69200
$ws.Range("E100").Value = -50100
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -8100
$ws.Range("E102").Value = 33400

# A handful of historical quarters (columns now H and I, and one cell in J)
# were also restated with corrected figures as part of this data refresh.
$ws.Range("I9").Value = 110400
$ws.Range("I10").Value = 179000
$ws.Range("I17").Value = 235600
$ws.Range("I18").Value = 53800
$ws.Range("I20").Value = 33300
$ws.Range("I21").Value = 92400
$ws.Range("I23").Value = 84400
$ws.Range("H24").Value = 21000
$ws.Range("I24").Value = 29500
$ws.Range("H26").Value = 35700
$ws.Range("I26").Value = 54900
$ws.Range("H27").Value = 35200
$ws.Range("I27").Value = 53600
$ws.Range("I32").Value = -33300
$ws.Range("H33").Value = 29800
$ws.Range("I33").Value = 53600
$ws.Range("H35").Value = 29800
$ws.Range("I35").Value = 53600
$ws.Range("H81").Value = 29800
$ws.Range("I81").Value = 53600
$ws.Range("H89").Value = 66000
$ws.Range("I89").Value = 17900
$ws.Range("I91").Value = -1700
$ws.Range("J91").Value = -1800
$ws.Range("H102").Value = -22500
$ws.Range("I102").Value = -204800
